$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 1: make it a bit taller (19.5pt, custom height) - column count (dimension/spans) will
# auto-extend to Q once column Q gets data below.
$ws.Rows.Item(1).RowHeight = 19.5

# Row 3: new bottom-border-only cell Q3 (same formatting as O3/P3).
$ws.Range("P3").Copy()
$ws.Range("Q3").PasteSpecial(-4122)

# Row 4: new year header cell Q4 = 2020 (same formatting as O4/P4).
$ws.Range("P4").Copy()
$ws.Range("Q4").PasteSpecial(-4122)
$ws.Range("Q4").Value = 2020

# Row 5: new data cell Q5 = 25.6
$ws.Range("P5").Copy()
$ws.Range("Q5").PasteSpecial(-4122)
$ws.Range("Q5").Value = 25.6

# Row 6: new data cell Q6
$ws.Range("P6").Copy()
$ws.Range("Q6").PasteSpecial(-4122)
$ws.Range("Q6").Value = 13.073527219449954

# Row 7: new data cell Q7
$ws.Range("P7").Copy()
$ws.Range("Q7").PasteSpecial(-4122)
$ws.Range("Q7").Value = 21.941290626870046

# Row 8: new data cell Q8
$ws.Range("P8").Copy()
$ws.Range("Q8").PasteSpecial(-4122)
$ws.Range("Q8").Value = 196.6
